# worked on export issue
#
# Territories format sheet: drop the CityGrade/IsActive columns (E:F),
# rename the RegionName header to CountryName and reuse the old CityName
# header cell as the new IsActive column, then insert two sample data
# rows (India/UP/Tamil Nadu and India/Bihar/Patna) under the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old CityGrade (E) and IsActive (F) columns entirely.
$ws.Range("E1:F1").EntireColumn.Delete()

# New column D take over as the IsActive flag column (was CityName).
$ws.Range("D1").Value = "IsActive"

# Make room for two sample rows right under the header.
$ws.Range("A2:A3").EntireRow.Insert()
$ws.Range("A2:D3").ClearFormats()

# Fill sample data (enter India before renaming A1, to mirror authoring order).
$ws.Range("A2").Value = "India"
$ws.Range("A1").Value = "CountryName"
$ws.Range("B2").Value = "UP"
$ws.Range("C2").Value = "Tamil Nadu"
$ws.Range("D2").Value = $true

$ws.Range("A3").Value = "India"
$ws.Range("B3").Value = "Bihar"
$ws.Range("C3").Value = "Patna"
$ws.Range("D3").Value = $true

# Re-create the True/False list validation on the (now) IsActive column D.
$dv = $ws.Range("D1:D1048576").Validation
$dv.Add(3, 1, 1, '"True,False"')
$dv.IgnoreBlank = $true
$dv.InCellDropdown = $true
$dv.ShowInput = $true
$dv.ShowError = $true

$ws.Range("D6").Select()
